# Apply the edits described by the commit diff:
#  - Shared string "Status" (header in F1) -> "STATUS"
#  - Active selection on Sheet1 changes from F6 to C6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Status" header cell (F1) to be upper-case "STATUS"
$ws.Range("F1").Value = "STATUS"

# Move the active selection/cell to C6 (was F6)
$ws.Range("C6").Select()
